$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (rows 2-11, columns A-I) as described by the commit "revisi 09 April 2025"
$data = @(
    @(45689.96666666667, 100,        2025, 2, 1,  23, 12, 32, 5),
    @(45691.35277777778, 69.61,      2025, 2, 3,  8,  28, 34, 0),
    @(45695.38680555556, 99.2,       2025, 2, 7,  9,  17, 38, 4),
    @(45695.80972222222, 77.41,      2025, 2, 7,  19, 26, 38, 4),
    @(45697.19375,        100,       2025, 2, 9,  4,  39, 40, 6),
    @(45698.33819444444, 100,        2025, 2, 10, 8,  7,  41, 0),
    @(45698.78472222222, 82.2,       2025, 2, 10, 18, 50, 41, 0),
    @(45699.42291666667, 96.84999999999999, 2025, 2, 11, 10, 9,  42, 1),
    @(45703.22569444445, 95.41,      2025, 2, 15, 5,  25, 46, 5),
    @(45708.31319444445, 97.06,      2025, 2, 20, 7,  31, 51, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
